$funds = @(
    @(0, "519133", "海富通改革驱动灵活配置混合", "175.25", "93.64", "2.03", "3.5576", 7),
    @(1, "159825", "富国中证农业主题ETF", "29.76", "99.63", "3.15", "0.9374", 10),
    @(2, "161017", "富国中证500指数增强(LOF)", "82.26", "91.34", "1.02", "0.8391", 3),
    @(3, "014420", "中欧成长领航一年持有混合A", "24.26", "50.50", "2.59", "0.6283", 7),
    @(4, "010080", "中欧优势成长三个月定期开放混合", "12.11", "85.82", "3.12", "0.3778", 6),
    @(5, "010790", "海富通均衡甄选混合A", "14.02", "93.79", "1.95", "0.2734", 8),
    @(6, "010770", "天弘中证农业主题指数C", "8.31", "95.31", "3.02", "0.2510", 10),
    @(7, "010769", "天弘中证农业主题指数A", "8.21", "95.31", "3.02", "0.2479", 10),
    @(8, "011429", "前海开源民裕进取混合", "2.83", "79.91", "7.17", "0.2029", 5),
    @(9, "159606", "易方达中证500质量成长ETF", "9.65", "98.03", "1.84", "0.1776", 7),
    @(10, "007593", "鹏扬中证500质量成长指数A", "9.59", "94.61", "1.76", "0.1688", 7),
    @(11, "014421", "中欧成长领航一年持有混合C", "4.88", "50.50", "2.59", "0.1264", 7),
    @(12, "010534", "广发均衡增长混合A", "30.75", "46.63", "0.33", "0.1015", 6),
    @(13, "010568", "海富通惠睿精选混合A", "15.75", "29.62", "0.60", "0.0945", 10),
    @(14, "011410", "中信建投量化进取6个月持有期混合A", "9.13", "93.80", "1.01", "0.0922", 5),
    @(15, "006048", "长城中证500指数增强A", "4.45", "92.64", "2.03", "0.0903", 8),
    @(16, "460009", "华泰柏瑞量化先行混合A", "9.13", "90.47", "0.89", "0.0813", 6),
    @(17, "010791", "海富通均衡甄选混合C", "4.08", "93.79", "1.95", "0.0796", 8),
    @(18, "007594", "鹏扬中证500质量成长指数C", "4.34", "94.61", "1.76", "0.0764", 7),
    @(19, "011588", "前海开源成份精选混合", "0.96", "81.61", "6.99", "0.0671", 4),
    @(20, "012080", "易方达中证500指数量化增强型证券投资基金A", "6.82", "84.83", "0.90", "0.0614", 8),
    @(21, "159827", "银华中证农业主题ETF", "1.83", "97.33", "3.10", "0.0567", 10),
    @(22, "014344", "鹏华中证500指数增强A", "1.99", "92.63", "1.86", "0.0370", 4),
    @(23, "007413", "长城中证500指数增强C", "1.72", "92.64", "2.03", "0.0349", 8),
    @(24, "010569", "海富通惠睿精选混合C", "5.79", "29.62", "0.60", "0.0347", 10),
    @(25, "560500", "鹏扬中证500质量成长交易型开放式指数证券投资基金", "1.71", "98.64", "1.84", "0.0315", 7),
    @(26, "005055", "华泰柏瑞量化阿尔法灵活配置混合A", "2.53", "89.49", "1.03", "0.0261", 8),
    @(27, "011411", "中信建投量化进取6个月持有期混合C", "2.15", "93.80", "1.01", "0.0217", 5),
    @(28, "003241", "创金合信量化发现灵活配置混合A", "2.07", "90.96", "1.05", "0.0217", 10),
    @(29, "010153", "中加中证500指数增强A", "1.26", "94.19", "1.63", "0.0205", 9),
    @(30, "014345", "鹏华中证500指数增强C", "0.78", "92.63", "1.86", "0.0145", 4),
    @(31, "012081", "易方达中证500指数量化增强型证券投资基金C", "1.57", "84.83", "0.90", "0.0141", 8),
    @(32, "006354", "国泰民裕进取灵活配置混合", "0.52", "80.18", "2.39", "0.0124", 1),
    @(33, "010154", "中加中证500指数增强C", "0.60", "94.19", "1.63", "0.0098", 9),
    @(34, "006783", "红土创新中证500指数增强A", "0.42", "91.83", "2.18", "0.0092", 8),
    @(35, "003242", "创金合信量化发现灵活配置混合C", "0.81", "90.96", "1.05", "0.0085", 10),
    @(36, "010535", "广发均衡增长混合C", "1.09", "46.63", "0.33", "0.0036", 6),
    @(37, "006784", "红土创新中证500指数增强C", "0.13", "91.83", "2.18", "0.0028", 8),
    @(38, "010246", "华泰柏瑞量化先行混合C", "0.12", "90.47", "0.89", "0.0011", 6),
    @(39, "006601", "国融融泰灵活配置混合A", "0.04", "47.44", "2.26", "0.0009", 10),
    @(40, "006602", "国融融泰灵活配置混合C", "0.01", "47.44", "2.26", "0.0002", 10),
    @(41, "006532", "华泰柏瑞量化阿尔法灵活配置混合C", "0.01", "89.49", "1.03", "0.0001", 8)
)

$wb = $excel.ActiveWorkbook

# --- Step 1: create the new "2022-Q1" sheet right before the "总计" sheet ---
$total = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($total)
$newSheet.Name = "2022-Q1"

# NOTE: inserting a sheet shifts worksheet positions, and `$total` was bound
# to the position "总计" used to occupy (now held by the new sheet) - so we
# must re-resolve it by name again before touching it.
$total = $wb.Worksheets.Item("总计")

# Header row
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# Make sure the text-like columns (B..G) are stored as text, not auto-converted numbers,
# so leading zeros in fund codes and trailing zeros in percentages survive.
$newSheet.Range("B2:G43").NumberFormat = "@"

foreach ($row in $funds) {
    $r = [int]$row[0] + 2
    $newSheet.Cells.Item($r,1).Value = $row[0]
    $newSheet.Cells.Item($r,2).Value = $row[1]
    $newSheet.Cells.Item($r,3).Value = $row[2]
    $newSheet.Cells.Item($r,4).Value = $row[3]
    $newSheet.Cells.Item($r,5).Value = $row[4]
    $newSheet.Cells.Item($r,6).Value = $row[5]
    $newSheet.Cells.Item($r,7).Value = $row[6]
    $newSheet.Cells.Item($r,8).Value = $row[7]
}

# Match the bold/centered/bordered header style used by the other sheets (style index 2),
# by copying formats from an existing sheet that already has it.
$styleSrc = $wb.Worksheets.Item("2021-Q4")
$styleSrc.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$styleSrc.Range("A2").Copy()
$newSheet.Range("A2:A43").PasteSpecial(-4122)

# --- Step 2: insert the new "2022-Q1" row at the top of the "总计" summary sheet ---
$total.Range("A2:D6").Copy($total.Range("A3:D7"))

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 42
$total.Cells.Item(2,4).Value = 8.890000000000001

for ($r = 3; $r -le 7; $r++) {
    $total.Cells.Item($r,1).Value = $r - 2
}
